# Journal de travail : ajout d'une nouvelle entrée (ligne 23) décrivant
# l'affichage de la grille et la vérification des coordonnées, faite au
# domicile le 2021-03-09 de 18h40 à 19h50, et mise à jour de la cellule
# sélectionnée.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# La ligne 23 du tableau (Tableau1) est encore vide : on récupère le format
# (format de date / heure) de la ligne précédente (22) pour les colonnes
# Date / Heure Début / Heure fin avant d'y inscrire les nouvelles valeurs.
$ws.Range("E22").Copy()
$ws.Range("E23").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F22").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G22").Copy()
$ws.Range("G23").PasteSpecial(-4122)

# Nouvelle entrée du journal : 09.03.2021, de 18:40 à 19:50.
# On renseigne "Heure fin" (G) avant "Heure Début" (F) afin que la formule
# de durée (colonne H, référence structurée du tableau) se recalcule
# correctement dès que la ligne, auparavant vide, reçoit sa dernière valeur.
$ws.Range("E23").Value = 44264
$ws.Range("G23").Value = 0.82638888888888884
$ws.Range("F23").Value = 0.77777777777777779

$ws.Range("I23").Value = "Développement"
$ws.Range("J23").Value = "Finalization de la grille, commencement des coordonnées"
$ws.Range("K23").Value = "Domicile"
$ws.Range("L23").Value = "Finalization de la grille et vérifiquation des coordonnées"

# La sélection active se déplace sur L24.
$ws.Range("L24").Select()
